$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from serial date 45185 to 45204 for rows 2-10
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value = 45204
}
